$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Start Frame value for the "Run" row (F5): 181 -> 210
$ws.Range("F5").Value = 210

# Move the active selection to F8 (matches the saved cursor position)
[void]$ws.Range("F8").Select()
